$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws "D2" "68.553.89"
Set-TextValue $ws "E2" "  +1.34%  "
Set-TextValue $ws "D3" "2.650.41"
Set-TextValue $ws "E3" "  +1.41%  "
Set-TextValue $ws "E4" "  -0.02%  "
Set-TextValue $ws "D5" "599.44"
Set-TextValue $ws "E5" "  +0.68%  "
Set-TextValue $ws "D6" "154.91"
Set-TextValue $ws "E6" "  +1.63%  "
Set-TextValue $ws "E7" "  -0.04%  "
Set-TextValue $ws "E8" "  +0.91%  "
Set-TextValue $ws "D9" "2.649.08"
Set-TextValue $ws "E9" "  +1.39%  "
Set-TextValue $ws "E10" "  +8.78%  "
Set-TextValue $ws "D11" "0.158"
Set-TextValue $ws "E11" "  -0.30%  "
Set-TextValue $ws "D12" "5.28"
Set-TextValue $ws "E12" "  +1.70%  "
Set-TextValue $ws "D13" "0.356"
Set-TextValue $ws "E13" "  +2.81%  "
Set-TextValue $ws "E14" "  +2.92%  "
Set-TextValue $ws "E15" "  +2.14%  "
Set-TextValue $ws "D16" "3.128.99"
Set-TextValue $ws "E16" "  +1.19%  "
Set-TextValue $ws "D17" "68.344.68"
Set-TextValue $ws "E17" "  +1.14%  "
Set-TextValue $ws "D18" "2.654.67"
Set-TextValue $ws "E18" "  +1.52%  "
Set-TextValue $ws "E19" "  +2.93%  "
Set-TextValue $ws "D20" "366.89"
Set-TextValue $ws "E20" "  -1.34%  "
Set-TextValue $ws "D21" "7.53"
Set-TextValue $ws "E21" "  +4.86%  "
Set-TextValue $ws "E22" "  +4.42%  "
Set-TextValue $ws "D24" "2.10"
Set-TextValue $ws "E24" "  +2.50%  "
Set-TextValue $ws "D25" "73.85"
Set-TextValue $ws "E25" "  +1.35%  "
Set-TextValue $ws "E26" "  +0.03%  "
Set-TextValue $ws "D27" "9.92"
Set-TextValue $ws "E27" "  +0.58%  "
Set-TextValue $ws "E28" "  +3.93%  "
Set-TextValue $ws "D29" "2.780.74"
Set-TextValue $ws "E29" "  +1.18%  "
Set-TextValue $ws "D30" "582.60"
Set-TextValue $ws "E30" "  -2.48%  "
Set-TextValue $ws "D31" "1.00"
Set-TextValue $ws "E31" "  -0.08%  "
Set-TextValue $ws "D32" "8.26"
Set-TextValue $ws "E32" "  +5.95%  "
Set-TextValue $ws "D33" "1.45"
Set-TextValue $ws "E33" "  +4.53%  "
Set-TextValue $ws "E34" "  +1.94%  "
Set-TextValue $ws "E35" "  +5.34%  "
Set-TextValue $ws "E36" "  +6.29%  "
Set-TextValue $ws "E37" "  +0.00%  "
Set-TextValue $ws "D38" "159.37"
Set-TextValue $ws "E38" "  +0.61%  "
Set-TextValue $ws "E39" "  +1.98%  "
Set-TextValue $ws "E40" "  +1.27%  "
Set-TextValue $ws "E41" "  +2.46%  "
Set-TextValue $ws "D42" "5.47"
Set-TextValue $ws "E42" "  +3.86%  "
Set-TextValue $ws "B43" "BabyDogeCoin"
Set-TextValue $ws "C43" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws "D43" "0.0₆0341"
Set-TextValue $ws "E43" "  +15.39%  "
Set-TextValue $ws "B44" "dogwifhat"
Set-TextValue $ws "C44" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws "D44" "2.71"
Set-TextValue $ws "E44" "  -0.38%  "
Set-TextValue $ws "E45" "  +3.62%  "
Set-TextValue $ws "E46" "  +0.02%  "
Set-TextValue $ws "D47" "40.49"
Set-TextValue $ws "E47" "  +0.28%  "
Set-TextValue $ws "D48" "158.09"
Set-TextValue $ws "E48" "  +1.31%  "
Set-TextValue $ws "E49" "  +3.42%  "
Set-TextValue $ws "E50" "  +2.56%  "
Set-TextValue $ws "D51" "22.11"
Set-TextValue $ws "E51" "  +4.18%  "
